$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 8: new time entry.
$ws.Range("A8").Value = "Udarbejdelse af SSD01 med Tommy"
$ws.Range("B8").Value = "Business-Process Analyst"
$ws.Range("C8").Value = 43886
$ws.Range("D8").Value = 0.35416666666666669
$ws.Range("E8").Value = 0.42708333333333331

# Row 9: new time entry.
$ws.Range("A9").Value = "Review ad AD01 med Matias"
$ws.Range("B9").Value = "Reviewer"
$ws.Range("C9").Value = 43886
$ws.Range("D9").Value = 0.4375
$ws.Range("E9").Value = 0.46875

# Row 3: fix capitalisation of the task description ("anders" -> "Anders").
$ws.Range("A3").Value = "Udarbejdelse af DOM08 med Anders"

# Row 10: new time entry.
$ws.Range("A10").Value = "Review af DOM01 med Toke"
$ws.Range("B10").Value = "Reviewer"
$ws.Range("C10").Value = 43886
$ws.Range("D10").Value = 0.50694444444444442
$ws.Range("E10").Value = 0.53819444444444442

# Row 11: new time entry.
$ws.Range("A11").Value = "Udarbejdelse af iterationsplan 2"
$ws.Range("B11").Value = "Project Manager"
$ws.Range("C11").Value = 43886
$ws.Range("D11").Value = 0.55208333333333337
$ws.Range("E11").Value = 0.70833333333333337

# Update the selected cell to match the saved view state.
$ws.Activate()
$ws.Range("D12").Select()
